$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Allergic reactions, including hives and anaphylaxis (in rare cases)" /
# "Secondary Infections:" / "Scratching the bites ..." / "Psychological Impact:"
# rows were mis-aligned (columns C/D shifted by one row). Split the combined
# "Allergic reactions..." text into its own Title/Description pair and shift
# the remaining bed_bug health-concern rows back into alignment.
$ws.Range("C7").Value = "Allergic reactions"
$ws.Range("D7").Value = "including hives and anaphylaxis (in rare cases)"

$ws.Range("C8").Value = "Secondary Infections:"
$ws.Range("D8").Value = "Scratching the bites can break the skin and increase the risk of bacterial or fungal infections."

# Match the "category title" look the other Health concern sub-headers use
# (C6 "Skin Irritation" already carries this alignment style).
$ws.Range("C7").WrapText = $false

# Update the saved selection to span the full data table instead of the
# stray D12 cell left over from editing.
$ws.Range("A1:D26").Select()
